$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (row 117), shifting
# the existing rows 117:145 down to 119:147. This represents a new week of
# reported prices (fecha = 44463) being added ahead of the historical rows.
$ws.Rows("117:118").Insert()

# New row 117 (Primera)
$ws.Range("A117").Value = 11
$ws.Range("B117").Value = "Vega Monumental Concepción"
$ws.Range("C117").Value = "Bíobío"
$ws.Range("D117").Value = 44463
$ws.Range("E117").Value = 8
$ws.Range("F117").Value = 100112008
$ws.Range("G117").Value = "Coliflor"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 2000
$ws.Range("K117").Value = 700
$ws.Range("L117").Value = 800
$ws.Range("M117").Value = 750
$ws.Range("N117").Value = "$/unidad"
$ws.Range("O117").Value = "Región Metropolitana"
$ws.Range("P117").Value = 750
$ws.Range("Q117").Value = 1
$ws.Range("R117").Value = "Hortaliza"

# New row 118 (Segunda)
$ws.Range("A118").Value = 11
$ws.Range("B118").Value = "Vega Monumental Concepción"
$ws.Range("C118").Value = "Bíobío"
$ws.Range("D118").Value = 44463
$ws.Range("E118").Value = 8
$ws.Range("F118").Value = 100112008
$ws.Range("G118").Value = "Coliflor"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Segunda"
$ws.Range("J118").Value = 1000
$ws.Range("K118").Value = 600
$ws.Range("L118").Value = 600
$ws.Range("M118").Value = 600
$ws.Range("N118").Value = "$/unidad"
$ws.Range("O118").Value = "Región Metropolitana"
$ws.Range("P118").Value = 600
$ws.Range("Q118").Value = 1
$ws.Range("R118").Value = "Hortaliza"

# Keep the D column's custom date number format (style index 2) consistent
# with the rest of the column for the two newly inserted rows.
$ws.Range("D117:D118").NumberFormat = $ws.Range("D119").NumberFormat
